$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..55 (54 records). A new weekly record is
# inserted, so a second new row is needed too: insert a fresh row at 41 which
# pushes the existing rows 41..55 down to 42..56, then fill rows 40 and 41
# with the two new records.
$ws.Rows(41).Insert()

# Row 40 - brand new record (previously occupied by the old row-40 record,
# which is now superseded by this one)
$ws.Cells.Item(40, 1).Value = 7
$ws.Cells.Item(40, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(40, 3).Value = "Ñuble"
$ws.Cells.Item(40, 4).Value = 44845
$ws.Cells.Item(40, 5).Value = 16
$ws.Cells.Item(40, 6).Value = 100112026
$ws.Cells.Item(40, 7).Value = "Haba"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 120
$ws.Cells.Item(40, 11).Value = 9000
$ws.Cells.Item(40, 12).Value = 9500
$ws.Cells.Item(40, 13).Value = 9250
$ws.Cells.Item(40, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 370
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"

# Row 41 - second new record (fills the row created by the insert above)
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = 44837
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112026
$ws.Cells.Item(41, 7).Value = "Haba"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 30
$ws.Cells.Item(41, 11).Value = 9000
$ws.Cells.Item(41, 12).Value = 9000
$ws.Cells.Item(41, 13).Value = 9000
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(41, 16).Value = 360
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
